## Mithamain project data collected completely
##
## The "main" sheet pulls its raw input numbers from an external workbook
## ("Hydrometer_V8.xlsm" -> 'input-output' sheet) that isn't available in
## this environment. Those source numbers were refreshed/updated; since the
## external workbook can't be reopened here, we push the refreshed values
## directly into the handful of "leaf" cells that previously mirrored the
## external cache. Every other figure on "main" and on "Report" (deltas,
## percentages, the trend/log numbers, and the chart feeding off Report)
## is a normal in-workbook formula, so it recalculates on its own once the
## leaf inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Row 16 block (C16:G16) -> input-output!D24:H24 -- only F24 (main!E16) moved.
$ws.Range("E16").Value2 = 25

# Row 17 block (C17:G17) -> input-output!D25:H25
$ws.Range("C17").Value2 = 7.33
$ws.Range("D17").Value2 = 10.83
$ws.Range("E17").Value2 = 10.83
$ws.Range("F17").Value2 = 7.35
$ws.Range("G17").Value2 = 9.05

# Row 18 block (C18:G18) -> input-output!D26:H26
$ws.Range("C18").Value2 = 43.41
$ws.Range("D18").Value2 = 42.58
$ws.Range("E18").Value2 = 41.74
$ws.Range("F18").Value2 = 41.36
$ws.Range("G18").Value2 = 41.09

# Row 19 block (C19:G19) -> input-output!D27:H27
$ws.Range("C19").Value2 = 32.02
$ws.Range("D19").Value2 = 32.84
$ws.Range("E19").Value2 = 33
$ws.Range("F19").Value2 = 31.26
$ws.Range("G19").Value2 = 31.72

# Row 44 block (D44:F44) -> input-output!M24:O24
$ws.Range("D44").Value2 = 7.1
$ws.Range("E44").Value2 = 10.23
$ws.Range("F44").Value2 = 9.8

# Row 45 block (D45:F45) -> input-output!M25:O25
$ws.Range("D45").Value2 = 41.97
$ws.Range("E45").Value2 = 46.2
$ws.Range("F45").Value2 = 43.17

# Row 46 block (D46:F46) -> input-output!M26:O26
$ws.Range("D46").Value2 = 35
$ws.Range("E46").Value2 = 38.83
$ws.Range("F46").Value2 = 36.39

# Everything downstream (C20:G22, D47:F49, G39:G41, G54, the Report sheet
# mirror formulas, and the chart's cached series) is formula-driven off the
# cells above, so a full recalc propagates the refreshed numbers.
$excel.CalculateFullRebuild()
